# Update gh-pages to output generated at 456a3b4
# Sets the "想去人数" (F column) counters to 0 for specific rows
# on the "展览" sheet (all three data rows with values) and on the
# "全部类型" sheet (rows 3 and 5 only; row 2 stays at 163).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 0
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F5").Value = 0

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 0
$wsAll.Range("F5").Value = 0
